$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I3").Value = 'qy'
$ws.Range("J3").Value = 'Yes-No-Question'
$ws.Range("I12").Value = 'sd'
$ws.Range("J12").Value = 'Statement-non-opinion'
$ws.Range("I15").Value = 'sv'
$ws.Range("J15").Value = 'Statement-opinion'
$ws.Range("I21").Value = 'b'
$ws.Range("J21").Value = 'Acknowledge (Backchannel)'
$ws.Range("I52").Value = 'sd'
$ws.Range("J52").Value = 'Statement-non-opinion'
$ws.Range("I81").Value = 'sv'
$ws.Range("J81").Value = 'Statement-opinion'
$ws.Range("I82").Value = 'aa'
$ws.Range("J82").Value = 'Agree/Accept'
$ws.Range("I85").Value = 'sv'
$ws.Range("J85").Value = 'Statement-opinion'
$ws.Range("I90").Value = '%'
$ws.Range("J90").Value = 'Uninterpretable'
$ws.Range("I94").Value = 'sd'
$ws.Range("J94").Value = 'Statement-non-opinion'
$ws.Range("I99").Value = 'sv'
$ws.Range("J99").Value = 'Statement-opinion'
$ws.Range("I101").Value = 'sv'
$ws.Range("J101").Value = 'Statement-opinion'
$ws.Range("I102").Value = 'sv'
$ws.Range("J102").Value = 'Statement-opinion'
$ws.Range("I104").Value = 'b'
$ws.Range("J104").Value = 'Acknowledge (Backchannel)'
$ws.Range("I116").Value = 'sv'
$ws.Range("J116").Value = 'Statement-opinion'
$ws.Range("I117").Value = 'aa'
$ws.Range("J117").Value = 'Agree/Accept'
$ws.Range("I118").Value = 'sv'
$ws.Range("J118").Value = 'Statement-opinion'
$ws.Range("I124").Value = 'b'
$ws.Range("J124").Value = 'Acknowledge (Backchannel)'
$ws.Range("I126").Value = 'b'
$ws.Range("J126").Value = 'Acknowledge (Backchannel)'
